$d = $word.ActiveDocument
$d.Content.Find.Execute("12/07", $true, $false, $false, $false, $false,
                         $true, 1, $false, "19/07", 2)
